$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2083333333333333
$ws.Range("C2").Value = 0.5208333333333334
$ws.Range("J2").Value = 0.02083333333333333
$ws.Range("P2").Value = 0.1547619047619048
$ws.Range("S2").Value = 0.09523809523809523
$ws.Range("B3").Value = 0.005494505494505495
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.02197802197802198
$ws.Range("P3").Value = 0.6813186813186813
$ws.Range("S3").Value = 0.2692307692307692
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.6818181818181818
$ws.Range("S4").Value = 0.2272727272727273
$ws.Range("B6").Value = 0.06172839506172839
$ws.Range("D6").Value = 0.01234567901234568
$ws.Range("F6").Value = 0.06172839506172839
$ws.Range("J6").Value = 0.2386831275720165
$ws.Range("O6").Value = 0.03292181069958848
$ws.Range("Q6").Value = 0.1563786008230453
$ws.Range("R6").Value = 0.06172839506172839
$ws.Range("S6").Value = 0.3744855967078189
$ws.Range("B7").Value = 0.08796296296296297
$ws.Range("D7").Value = 0.01388888888888889
$ws.Range("F7").Value = 0.07870370370370371
$ws.Range("J7").Value = 0.1064814814814815
$ws.Range("O7").Value = 0.01388888888888889
$ws.Range("Q7").Value = 0.1481481481481481
$ws.Range("R7").Value = 0.07407407407407407
$ws.Range("S7").Value = 0.4768518518518519
$ws.Range("B8").Value = 0.1076233183856502
$ws.Range("D8").Value = 0.01121076233183856
$ws.Range("E8").Value = 0.002242152466367713
$ws.Range("F8").Value = 0.05605381165919283
$ws.Range("J8").Value = 0.1008968609865471
$ws.Range("O8").Value = 0.0336322869955157
$ws.Range("Q8").Value = 0.1390134529147982
$ws.Range("R8").Value = 0.1345291479820628
$ws.Range("S8").Value = 0.4147982062780269
$ws.Range("B9").Value = 0.1126760563380282
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.06103286384976526
$ws.Range("J9").Value = 0.1220657276995305
$ws.Range("O9").Value = 0.0187793427230047
$ws.Range("Q9").Value = 0.1784037558685446
$ws.Range("R9").Value = 0.05633802816901409
$ws.Range("S9").Value = 0.4366197183098591
$ws.Range("B10").Value = 0.1139896373056995
$ws.Range("D10").Value = 0.0229459659511473
$ws.Range("E10").Value = 0.003700962250185048
$ws.Range("F10").Value = 0.0695780903034789
$ws.Range("J10").Value = 0.1028867505551443
$ws.Range("O10").Value = 0.02146558105107328
$ws.Range("Q10").Value = 0.2102146558105107
$ws.Range("R10").Value = 0.07105847520355292
$ws.Range("S10").Value = 0.384159881569208
$ws.Range("F11").Value = 0.003048780487804878
$ws.Range("G11").Value = 0.1524390243902439
$ws.Range("J11").Value = 0.08231707317073171
$ws.Range("K11").Value = 0.1920731707317073
$ws.Range("L11").Value = 0.5579268292682927
$ws.Range("S11").Value = 0.01219512195121951
$ws.Range("G12").Value = 0.7553191489361702
$ws.Range("J12").Value = 0.1968085106382979
$ws.Range("K12").Value = 0.005319148936170213
$ws.Range("L12").Value = 0.01595744680851064
$ws.Range("S12").Value = 0.02659574468085106
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2790697674418605
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.08583690987124463
$ws.Range("I15").Value = 0.03433476394849785
$ws.Range("J15").Value = 0.4034334763948498
$ws.Range("K15").Value = 0.06437768240343347
$ws.Range("M15").Value = 0.0128755364806867
$ws.Range("O15").Value = 0.0815450643776824
$ws.Range("S15").Value = 0.3004291845493562
$ws.Range("F16").Value = 0.02926829268292683
$ws.Range("H16").Value = 0.1707317073170732
$ws.Range("I16").Value = 0.1073170731707317
$ws.Range("J16").Value = 0.4097560975609756
$ws.Range("K16").Value = 0.1073170731707317
$ws.Range("M16").Value = 0.01951219512195122
$ws.Range("O16").Value = 0.03414634146341464
$ws.Range("S16").Value = 0.1219512195121951
$ws.Range("F17").Value = 0.01106194690265487
$ws.Range("H17").Value = 0.1592920353982301
$ws.Range("I17").Value = 0.07964601769911504
$ws.Range("J17").Value = 0.4358407079646018
$ws.Range("K17").Value = 0.1039823008849557
$ws.Range("M17").Value = 0.01106194690265487
$ws.Range("O17").Value = 0.05309734513274336
$ws.Range("S17").Value = 0.1460176991150443
$ws.Range("F18").Value = 0.03015075376884422
$ws.Range("H18").Value = 0.08542713567839195
$ws.Range("I18").Value = 0.1005025125628141
$ws.Range("J18").Value = 0.407035175879397
$ws.Range("K18").Value = 0.1658291457286432
$ws.Range("M18").Value = 0.01507537688442211
$ws.Range("O18").Value = 0.05527638190954774
$ws.Range("S18").Value = 0.1407035175879397
$ws.Range("F19").Value = 0.01931034482758621
$ws.Range("H19").Value = 0.2110344827586207
$ws.Range("I19").Value = 0.0896551724137931
$ws.Range("J19").Value = 0.3675862068965517
$ws.Range("K19").Value = 0.1020689655172414
$ws.Range("M19").Value = 0.02068965517241379
$ws.Range("N19").Value = 0.001379310344827586
$ws.Range("O19").Value = 0.05931034482758621
$ws.Range("S19").Value = 0.1289655172413793
